$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows (9:10) for rho_m / rho_o, shifting everything below down by 2 ---
$ws.Rows("9:10").Insert()

$ws.Range("A9").Value = "rho_m"
$ws.Range("B9").Value = 2650
$ws.Range("D9").Value = "density quartz"

$ws.Range("A10").Value = "rho_o"
$ws.Range("B10").Value = 2650
$ws.Range("D10").Value = "not sure"

# --- Rebuild the Van Genuchten parameter block (now rows 38-43 after the shift) ---
# Original block held: alpha_sand, alpha_silt, n_sand, n_silt, residual_wc_sand,
# residual_wc_silt. New block adds water/clay/peat rows and expands to rows 38-52.
$ws.Rows("38:43").ClearContents()

$ws.Range("A38").Value = "alpha_water"
$ws.Range("B38").Value = 400
$ws.Range("B38").NumberFormat = "0.00"

$ws.Range("A39").Value = "alpha_sand"
$ws.Range("B39").Value = 4.06
$ws.Range("B39").NumberFormat = "0.00"

$ws.Range("A40").Value = "alpha_silt"
$ws.Range("B40").Value = 0.65
$ws.Range("B40").NumberFormat = "0.00"

$ws.Range("A41").Value = "alpha_clay"
$ws.Range("B41").Value = 1.49
$ws.Range("B41").NumberFormat = "0.00"

$ws.Range("A42").Value = "alpha_peat"
$ws.Range("B42").Value = 2.31
$ws.Range("B42").NumberFormat = "0.00"
$ws.Range("D42").Value = "from Hydraulic properties of fen peat soils in Poland, Gnatowski 2010"

$ws.Range("A43").Value = "n_water"
$ws.Range("B43").Value = 2.5
$ws.Range("B43").NumberFormat = "0.00"

$ws.Range("A44").Value = "n_sand"
$ws.Range("B44").Value = 2
$ws.Range("B44").NumberFormat = "0.00"

$ws.Range("A45").Value = "n_silt"
$ws.Range("B45").Value = 1.7
$ws.Range("B45").NumberFormat = "0.00"

$ws.Range("A46").Value = "n_clay"
$ws.Range("B46").Value = 1.25
$ws.Range("B46").NumberFormat = "0.00"

$ws.Range("A47").Value = "n_peat"
$ws.Range("B47").Value = 1.292
$ws.Range("B47").NumberFormat = "0.00"

$ws.Range("A48").Value = "residual_wc_water"
$ws.Range("B48").Value = 0
$ws.Range("B48").NumberFormat = "0.00"

$ws.Range("A49").Value = "residual_wc_sand"
$ws.Range("B49").Value = 0
$ws.Range("B49").NumberFormat = "0.00"

$ws.Range("A50").Value = "residual_wc_silt"
$ws.Range("B50").Value = 0
$ws.Range("B50").NumberFormat = "0.00"

$ws.Range("A51").Value = "residual_wc_clay"
$ws.Range("B51").Value = 0
$ws.Range("B51").NumberFormat = "0.00"

$ws.Range("A52").Value = "residual_wc_peat"
$ws.Range("B52").Value = 0
$ws.Range("B52").NumberFormat = "0.00"

# --- View state: scroll window & selection to match the saved workbook state ---
$ws.Range("D46").Select()
$excel.ActiveWindow.ScrollRow = 19
